$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1): insert a new "company" header in B1, shift the
# rest of the header labels right, and extend with the common
# property/legislator metadata headers used by the other sheets ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Copy the existing bold/bordered header style onto the newly added
# header cells (F1:K1) so they match B1:E1. (E1 already existed in the
# sheet with the correct header style, so it needs no copy/paste.)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Dates are stored as plain text ("2012-04-30"), not real date
# serials, in this workbook - format the date column as Text first so
# Excel does not auto-convert the literal into a date value.
$ws.Range("G2:G4").NumberFormat = "@"

# --- Row 2 (record 100): 國泰人壽 / 鍾愛一生313 ---
$ws.Range("A2").Value = 100
$ws.Range("B2").Value = "國泰人壽"
$ws.Range("C2").Value = "鍾愛一生313"
$ws.Range("D2").Value = "黃靜秋"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2012-04-30"
$ws.Range("H2").Value = "羅明才"
$ws.Range("I2").Value = 879
$ws.Range("J2").Value = "tmpa5201"
$ws.Range("K2").Value = 100

# --- Row 3 (record 101): 保德信國際人壽 / 教育終身壽險 ---
$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "保德信國際人壽"
$ws.Range("C3").Value = "教育終身壽險"
$ws.Range("D3").Value = "黃靜秋"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2012-04-30"
$ws.Range("H3").Value = "羅明才"
$ws.Range("I3").Value = 879
$ws.Range("J3").Value = "tmpa5201"
$ws.Range("K3").Value = 101

# --- Row 4 (record 102): 保德信國際人壽 / 教育終身壽險 ---
$ws.Range("A4").Value = 102
$ws.Range("B4").Value = "保德信國際人壽"
$ws.Range("C4").Value = "教育終身壽險"
$ws.Range("D4").Value = "黄靜秋"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("G4").Value = "2012-04-30"
$ws.Range("H4").Value = "羅明才"
$ws.Range("I4").Value = 879
$ws.Range("J4").Value = "tmpa5201"
$ws.Range("K4").Value = 102

# Copy the existing plain data-row style onto the newly added data
# cells (F2:K4) so they match B2:E4. (Column E already existed in the
# sheet with the correct data style, so it needs no copy/paste.)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("F2:K4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
